# Auto-generated Excel COM-interop script to refresh market-price data
# columns H-N (currentAveragePrice*, LevePrice*, LeveProfit*) for specific
# leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, matching a
# scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64 (ALC)
$ws.Cells.Item(64, 8).Value = 3251.4119  # H64: 3240.7368 -> 3251.4119
$ws.Cells.Item(64, 10).Value = 3297.3845  # J64: 3277.7334 -> 3297.3845
$ws.Cells.Item(64, 12).Value = 3297.3845  # L64: 3277.7334 -> 3297.3845
$ws.Cells.Item(64, 14).Value = -3793.3845  # N64: -3773.7334 -> -3793.3845

# Row 67 (ALC)
$ws.Cells.Item(67, 8).Value = 3251.4119  # H67: 3240.7368 -> 3251.4119
$ws.Cells.Item(67, 10).Value = 3297.3845  # J67: 3277.7334 -> 3297.3845
$ws.Cells.Item(67, 12).Value = 3297.3845  # L67: 3277.7334 -> 3297.3845
$ws.Cells.Item(67, 14).Value = -5013.3845  # N67: -4993.7334 -> -5013.3845

# Row 69 (ALC)
$ws.Cells.Item(69, 8).Value = 19166.666  # H69: 8871.200000000001 -> 19166.666
$ws.Cells.Item(69, 9).Value = 3000  # I69: 3013 -> 3000
$ws.Cells.Item(69, 10).Value = 27250  # J69: 9522.111000000001 -> 27250
$ws.Cells.Item(69, 11).Value = 9000  # K69: 9039 -> 9000
$ws.Cells.Item(69, 12).Value = 81750  # L69: 28566.333 -> 81750
$ws.Cells.Item(69, 13).Value = -8126  # M69: -8165 -> -8126
$ws.Cells.Item(69, 14).Value = -83498  # N69: -30314.333 -> -83498

# Row 72 (ALC)
$ws.Cells.Item(72, 8).Value = 19166.666  # H72: 8871.200000000001 -> 19166.666
$ws.Cells.Item(72, 9).Value = 3000  # I72: 3013 -> 3000
$ws.Cells.Item(72, 10).Value = 27250  # J72: 9522.111000000001 -> 27250
$ws.Cells.Item(72, 11).Value = 27000  # K72: 27117 -> 27000
$ws.Cells.Item(72, 12).Value = 245250  # L72: 85698.99900000001 -> 245250
$ws.Cells.Item(72, 13).Value = -22632  # M72: -22749 -> -22632
$ws.Cells.Item(72, 14).Value = -253986  # N72: -94434.99900000001 -> -253986

# Row 74 (ALC)
$ws.Cells.Item(74, 8).Value = 2761.2222  # H74: 3123.1765 -> 2761.2222
$ws.Cells.Item(74, 9).Value = 2607.6924  # I74: 2984.2856 -> 2607.6924
$ws.Cells.Item(74, 10).Value = 3160.4  # J74: 3220.4 -> 3160.4
$ws.Cells.Item(74, 11).Value = 2607.6924  # K74: 2984.2856 -> 2607.6924
$ws.Cells.Item(74, 12).Value = 3160.4  # L74: 3220.4 -> 3160.4
$ws.Cells.Item(74, 13).Value = -1671.6924  # M74: -2048.2856 -> -1671.6924
$ws.Cells.Item(74, 14).Value = -5032.4  # N74: -5092.4 -> -5032.4

# Row 76 (ALC)
$ws.Cells.Item(76, 8).Value = 3164.7058  # H76: 2696.6875 -> 3164.7058
$ws.Cells.Item(76, 9).Value = 3200  # I76: 2814.4285 -> 3200
$ws.Cells.Item(76, 10).Value = 3080  # J76: 2471.9092 -> 3080
$ws.Cells.Item(76, 11).Value = 3200  # K76: 2814.4285 -> 3200
$ws.Cells.Item(76, 12).Value = 3080  # L76: 2471.9092 -> 3080
$ws.Cells.Item(76, 13).Value = -2885  # M76: -2499.4285 -> -2885
$ws.Cells.Item(76, 14).Value = -3710  # N76: -3101.9092 -> -3710

# Row 77 (ALC)
$ws.Cells.Item(77, 8).Value = 2761.2222  # H77: 3123.1765 -> 2761.2222
$ws.Cells.Item(77, 9).Value = 2607.6924  # I77: 2984.2856 -> 2607.6924
$ws.Cells.Item(77, 10).Value = 3160.4  # J77: 3220.4 -> 3160.4
$ws.Cells.Item(77, 11).Value = 13038.462  # K77: 14921.428 -> 13038.462
$ws.Cells.Item(77, 12).Value = 15802  # L77: 16102 -> 15802
$ws.Cells.Item(77, 13).Value = -8358.462  # M77: -10241.428 -> -8358.462
$ws.Cells.Item(77, 14).Value = -25162  # N77: -25462 -> -25162

# Row 79 (ALC)
$ws.Cells.Item(79, 8).Value = 3164.7058  # H79: 2696.6875 -> 3164.7058
$ws.Cells.Item(79, 9).Value = 3200  # I79: 2814.4285 -> 3200
$ws.Cells.Item(79, 10).Value = 3080  # J79: 2471.9092 -> 3080
$ws.Cells.Item(79, 11).Value = 3200  # K79: 2814.4285 -> 3200
$ws.Cells.Item(79, 12).Value = 3080  # L79: 2471.9092 -> 3080
$ws.Cells.Item(79, 13).Value = -2108  # M79: -1722.4285 -> -2108
$ws.Cells.Item(79, 14).Value = -5264  # N79: -4655.9092 -> -5264

# Row 80 (ALC)
$ws.Cells.Item(80, 8).Value = 2370.95  # H80: 2602.0527 -> 2370.95
$ws.Cells.Item(80, 9).Value = 1958.8572  # I80: 2089 -> 1958.8572
$ws.Cells.Item(80, 10).Value = 2592.8462  # J80: 2975.182 -> 2592.8462
$ws.Cells.Item(80, 11).Value = 5876.571599999999  # K80: 6267 -> 5876.571599999999
$ws.Cells.Item(80, 12).Value = 7778.5386  # L80: 8925.545999999998 -> 7778.5386
$ws.Cells.Item(80, 13).Value = -4878.571599999999  # M80: -5269 -> -4878.571599999999
$ws.Cells.Item(80, 14).Value = -9774.5386  # N80: -10921.546 -> -9774.5386

# Row 83 (ALC)
$ws.Cells.Item(83, 8).Value = 2370.95  # H83: 2602.0527 -> 2370.95
$ws.Cells.Item(83, 9).Value = 1958.8572  # I83: 2089 -> 1958.8572
$ws.Cells.Item(83, 10).Value = 2592.8462  # J83: 2975.182 -> 2592.8462
$ws.Cells.Item(83, 11).Value = 17629.7148  # K83: 18801 -> 17629.7148
$ws.Cells.Item(83, 12).Value = 23335.6158  # L83: 26776.638 -> 23335.6158
$ws.Cells.Item(83, 13).Value = -12637.7148  # M83: -13809 -> -12637.7148
$ws.Cells.Item(83, 14).Value = -33319.6158  # N83: -36760.638 -> -33319.6158

# Row 101 (ALC)
$ws.Cells.Item(101, 8).Value = 1015.0909  # H101: 625.93335 -> 1015.0909
$ws.Cells.Item(101, 9).Value = 825.1429000000001  # I101: 788 -> 825.1429000000001
$ws.Cells.Item(101, 10).Value = 1347.5  # J101: 440.7143 -> 1347.5
$ws.Cells.Item(101, 11).Value = 2475.4287  # K101: 2364 -> 2475.4287
$ws.Cells.Item(101, 12).Value = 4042.5  # L101: 1322.1429 -> 4042.5
$ws.Cells.Item(101, 13).Value = -853.4287000000004  # M101: -742 -> -853.4287000000004
$ws.Cells.Item(101, 14).Value = -7286.5  # N101: -4566.1429 -> -7286.5

# Row 137 (ALC)
$ws.Cells.Item(137, 8).Value = 1381.2391  # H137: 1568.525 -> 1381.2391
$ws.Cells.Item(137, 9).Value = 1191.6786  # I137: 1480.5 -> 1191.6786
$ws.Cells.Item(137, 11).Value = 3575.0358  # K137: 4441.5 -> 3575.0358
$ws.Cells.Item(137, 13).Value = -1025.0358  # M137: -1891.5 -> -1025.0358

$ws = $wb.Worksheets.Item("ARM")
# Row 63 (ARM)
$ws.Cells.Item(63, 8).Value = 3760  # H63: 3992 -> 3760
$ws.Cells.Item(63, 9).Value = 2646.6667  # I63: 2485 -> 2646.6667
$ws.Cells.Item(63, 10).Value = 4714.2856  # J63: 5714.2856 -> 4714.2856
$ws.Cells.Item(63, 11).Value = 2646.6667  # K63: 2485 -> 2646.6667
$ws.Cells.Item(63, 12).Value = 4714.2856  # L63: 5714.2856 -> 4714.2856
$ws.Cells.Item(63, 13).Value = -1960.6667  # M63: -1799 -> -1960.6667
$ws.Cells.Item(63, 14).Value = -6086.2856  # N63: -7086.2856 -> -6086.2856

# Row 66 (ARM)
$ws.Cells.Item(66, 8).Value = 3760  # H66: 3992 -> 3760
$ws.Cells.Item(66, 9).Value = 2646.6667  # I66: 2485 -> 2646.6667
$ws.Cells.Item(66, 10).Value = 4714.2856  # J66: 5714.2856 -> 4714.2856
$ws.Cells.Item(66, 11).Value = 13233.3335  # K66: 12425 -> 13233.3335
$ws.Cells.Item(66, 12).Value = 23571.428  # L66: 28571.428 -> 23571.428
$ws.Cells.Item(66, 13).Value = -9801.333500000001  # M66: -8993 -> -9801.333500000001
$ws.Cells.Item(66, 14).Value = -30435.428  # N66: -35435.428 -> -30435.428

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (BSM)
$ws.Cells.Item(22, 8).Value = 598  # H22: 638.8 -> 598
$ws.Cells.Item(22, 9).Value = 598  # I22: 638.8 -> 598
$ws.Cells.Item(22, 11).Value = 598  # K22: 638.8 -> 598
$ws.Cells.Item(22, 13).Value = -425  # M22: -465.8 -> -425

$ws = $wb.Worksheets.Item("CRP")
# Row 64 (CRP)
$ws.Cells.Item(64, 8).Value = 27103.334  # H64: 27986.5 -> 27103.334
$ws.Cells.Item(64, 10).Value = 27103.334  # J64: 27986.5 -> 27103.334
$ws.Cells.Item(64, 12).Value = 27103.334  # L64: 27986.5 -> 27103.334
$ws.Cells.Item(64, 14).Value = -27599.334  # N64: -28482.5 -> -27599.334

# Row 67 (CRP)
$ws.Cells.Item(67, 8).Value = 27103.334  # H67: 27986.5 -> 27103.334
$ws.Cells.Item(67, 10).Value = 27103.334  # J67: 27986.5 -> 27103.334
$ws.Cells.Item(67, 12).Value = 27103.334  # L67: 27986.5 -> 27103.334
$ws.Cells.Item(67, 14).Value = -28819.334  # N67: -29702.5 -> -28819.334

$ws = $wb.Worksheets.Item("CUL")
# Row 40 (CUL)
$ws.Cells.Item(40, 8).Value = 301.0435  # H40: 290.13043 -> 301.0435
$ws.Cells.Item(40, 9).Value = 84.09090999999999  # I40: 81.083336 -> 84.09090999999999
$ws.Cells.Item(40, 10).Value = 499.91666  # J40: 518.1818 -> 499.91666
$ws.Cells.Item(40, 11).Value = 336.36364  # K40: 324.333344 -> 336.36364
$ws.Cells.Item(40, 12).Value = 1999.66664  # L40: 2072.7272 -> 1999.66664
$ws.Cells.Item(40, 13).Value = -267.36364  # M40: -255.333344 -> -267.36364
$ws.Cells.Item(40, 14).Value = -2137.66664  # N40: -2210.7272 -> -2137.66664

# Row 62 (CUL)
$ws.Cells.Item(62, 8).Value = 4037.6667  # H62: 3960.8572 -> 4037.6667
$ws.Cells.Item(62, 9).Value = 0  # I62: 2000 -> 0
$ws.Cells.Item(62, 10).Value = 4037.6667  # J62: 4111.6924 -> 4037.6667
$ws.Cells.Item(62, 11).Value = 0  # K62: 6000 -> 0
$ws.Cells.Item(62, 12).Value = 12113.0001  # L62: 12335.0772 -> 12113.0001
$ws.Cells.Item(62, 13).ClearContents()  # M62: -5314 -> (removed)
$ws.Cells.Item(62, 14).Value = -13485.0001  # N62: -13707.0772 -> -13485.0001

# Row 65 (CUL)
$ws.Cells.Item(65, 8).Value = 4037.6667  # H65: 3960.8572 -> 4037.6667
$ws.Cells.Item(65, 9).Value = 0  # I65: 2000 -> 0
$ws.Cells.Item(65, 10).Value = 4037.6667  # J65: 4111.6924 -> 4037.6667
$ws.Cells.Item(65, 11).Value = 0  # K65: 18000 -> 0
$ws.Cells.Item(65, 12).Value = 36339.0003  # L65: 37005.2316 -> 36339.0003
$ws.Cells.Item(65, 13).ClearContents()  # M65: -14568 -> (removed)
$ws.Cells.Item(65, 14).Value = -43203.0003  # N65: -43869.2316 -> -43203.0003

# Row 68 (CUL)
$ws.Cells.Item(68, 8).Value = 1679.8  # H68: 1314.1428 -> 1679.8
$ws.Cells.Item(68, 9).Value = 600  # I68: 450 -> 600
$ws.Cells.Item(68, 10).Value = 1949.75  # J68: 1659.8 -> 1949.75
$ws.Cells.Item(68, 11).Value = 1800  # K68: 1350 -> 1800
$ws.Cells.Item(68, 12).Value = 5849.25  # L68: 4979.4 -> 5849.25
$ws.Cells.Item(68, 13).Value = -989  # M68: -539 -> -989
$ws.Cells.Item(68, 14).Value = -7471.25  # N68: -6601.4 -> -7471.25

# Row 69 (CUL)
$ws.Cells.Item(69, 8).Value = 3720.7273  # H69: 5994 -> 3720.7273
$ws.Cells.Item(69, 9).Value = 1000  # I69: 0 -> 1000
$ws.Cells.Item(69, 10).Value = 4741  # J69: 5994 -> 4741
$ws.Cells.Item(69, 11).Value = 3000  # K69: 0 -> 3000
$ws.Cells.Item(69, 12).Value = 14223  # L69: 17982 -> 14223
$ws.Cells.Item(69, 13).Value = -2189  # M69: (new) -> -2189
$ws.Cells.Item(69, 14).Value = -15845  # N69: -19604 -> -15845

# Row 71 (CUL)
$ws.Cells.Item(71, 8).Value = 1679.8  # H71: 1314.1428 -> 1679.8
$ws.Cells.Item(71, 9).Value = 600  # I71: 450 -> 600
$ws.Cells.Item(71, 10).Value = 1949.75  # J71: 1659.8 -> 1949.75
$ws.Cells.Item(71, 11).Value = 5400  # K71: 4050 -> 5400
$ws.Cells.Item(71, 12).Value = 17547.75  # L71: 14938.2 -> 17547.75
$ws.Cells.Item(71, 13).Value = -1344  # M71: 6 -> -1344
$ws.Cells.Item(71, 14).Value = -25659.75  # N71: -23050.2 -> -25659.75

# Row 72 (CUL)
$ws.Cells.Item(72, 8).Value = 3720.7273  # H72: 5994 -> 3720.7273
$ws.Cells.Item(72, 9).Value = 1000  # I72: 0 -> 1000
$ws.Cells.Item(72, 10).Value = 4741  # J72: 5994 -> 4741
$ws.Cells.Item(72, 11).Value = 9000  # K72: 0 -> 9000
$ws.Cells.Item(72, 12).Value = 42669  # L72: 53946 -> 42669
$ws.Cells.Item(72, 13).Value = -4944  # M72: (new) -> -4944
$ws.Cells.Item(72, 14).Value = -50781  # N72: -62058 -> -50781

# Row 74 (CUL)
$ws.Cells.Item(74, 8).Value = 8009  # H74: 8012.6665 -> 8009
$ws.Cells.Item(74, 10).Value = 8009  # J74: 8012.6665 -> 8009
$ws.Cells.Item(74, 12).Value = 24027  # L74: 24037.9995 -> 24027
$ws.Cells.Item(74, 14).Value = -26149  # N74: -26159.9995 -> -26149

# Row 77 (CUL)
$ws.Cells.Item(77, 8).Value = 8009  # H77: 8012.6665 -> 8009
$ws.Cells.Item(77, 10).Value = 8009  # J77: 8012.6665 -> 8009
$ws.Cells.Item(77, 12).Value = 72081  # L77: 72113.9985 -> 72081
$ws.Cells.Item(77, 14).Value = -82689  # N77: -82721.9985 -> -82689

$ws = $wb.Worksheets.Item("GSM")
# Row 64 (GSM)
$ws.Cells.Item(64, 8).Value = 28163  # H64: 30000 -> 28163
$ws.Cells.Item(64, 10).Value = 28163  # J64: 30000 -> 28163
$ws.Cells.Item(64, 12).Value = 28163  # L64: 30000 -> 28163
$ws.Cells.Item(64, 14).Value = -28659  # N64: -30496 -> -28659

# Row 67 (GSM)
$ws.Cells.Item(67, 8).Value = 28163  # H67: 30000 -> 28163
$ws.Cells.Item(67, 10).Value = 28163  # J67: 30000 -> 28163
$ws.Cells.Item(67, 12).Value = 28163  # L67: 30000 -> 28163
$ws.Cells.Item(67, 14).Value = -29879  # N67: -31716 -> -29879

# Row 80 (GSM)
$ws.Cells.Item(80, 8).Value = 2295  # H80: 2290.9524 -> 2295
$ws.Cells.Item(80, 9).Value = 2356  # I80: 2363.3333 -> 2356
$ws.Cells.Item(80, 10).Value = 2218.75  # J80: 2194.4443 -> 2218.75
$ws.Cells.Item(80, 11).Value = 2356  # K80: 2363.3333 -> 2356
$ws.Cells.Item(80, 12).Value = 2218.75  # L80: 2194.4443 -> 2218.75
$ws.Cells.Item(80, 13).Value = -1358  # M80: -1365.3333 -> -1358
$ws.Cells.Item(80, 14).Value = -4214.75  # N80: -4190.4443 -> -4214.75

# Row 83 (GSM)
$ws.Cells.Item(83, 8).Value = 2295  # H83: 2290.9524 -> 2295
$ws.Cells.Item(83, 9).Value = 2356  # I83: 2363.3333 -> 2356
$ws.Cells.Item(83, 10).Value = 2218.75  # J83: 2194.4443 -> 2218.75
$ws.Cells.Item(83, 11).Value = 11780  # K83: 11816.6665 -> 11780
$ws.Cells.Item(83, 12).Value = 11093.75  # L83: 10972.2215 -> 11093.75
$ws.Cells.Item(83, 13).Value = -6788  # M83: -6824.666499999999 -> -6788
$ws.Cells.Item(83, 14).Value = -21077.75  # N83: -20956.2215 -> -21077.75

# Row 132 (GSM)
$ws.Cells.Item(132, 8).Value = 2634545.2  # H132: 3033516.8 -> 2634545.2
$ws.Cells.Item(132, 9).Value = 2672.1924  # I132: 3009.8262 -> 2672.1924
$ws.Cells.Item(132, 10).Value = 8336936.5  # J132: 10003683 -> 8336936.5
$ws.Cells.Item(132, 11).Value = 8016.5772  # K132: 9029.4786 -> 8016.5772
$ws.Cells.Item(132, 12).Value = 25010809.5  # L132: 30011049 -> 25010809.5
$ws.Cells.Item(132, 13).Value = -5486.5772  # M132: -6499.4786 -> -5486.5772
$ws.Cells.Item(132, 14).Value = -25015869.5  # N132: -30016109 -> -25015869.5

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (LTW)
$ws.Cells.Item(40, 8).Value = 2132.5264  # H40: 1853.8214 -> 2132.5264
$ws.Cells.Item(40, 9).Value = 2059.2942  # I40: 1853.8214 -> 2059.2942
$ws.Cells.Item(40, 10).Value = 2755  # J40: 0 -> 2755
$ws.Cells.Item(40, 11).Value = 2059.2942  # K40: 1853.8214 -> 2059.2942
$ws.Cells.Item(40, 12).Value = 2755  # L40: 0 -> 2755
$ws.Cells.Item(40, 13).Value = -1923.2942  # M40: -1717.8214 -> -1923.2942
$ws.Cells.Item(40, 14).Value = -3027  # N40: (new) -> -3027

# Row 82 (LTW)
$ws.Cells.Item(82, 8).Value = 1142.4546  # H82: 1071.4482 -> 1142.4546
$ws.Cells.Item(82, 9).Value = 853.3333  # I82: 856 -> 853.3333
$ws.Cells.Item(82, 10).Value = 1188.1052  # J82: 1153.5238 -> 1188.1052
$ws.Cells.Item(82, 11).Value = 853.3333  # K82: 856 -> 853.3333
$ws.Cells.Item(82, 12).Value = 1188.1052  # L82: 1153.5238 -> 1188.1052
$ws.Cells.Item(82, 13).Value = -492.3333  # M82: -495 -> -492.3333
$ws.Cells.Item(82, 14).Value = -1910.1052  # N82: -1875.5238 -> -1910.1052

# Row 85 (LTW)
$ws.Cells.Item(85, 8).Value = 1142.4546  # H85: 1071.4482 -> 1142.4546
$ws.Cells.Item(85, 9).Value = 853.3333  # I85: 856 -> 853.3333
$ws.Cells.Item(85, 10).Value = 1188.1052  # J85: 1153.5238 -> 1188.1052
$ws.Cells.Item(85, 11).Value = 853.3333  # K85: 856 -> 853.3333
$ws.Cells.Item(85, 12).Value = 1188.1052  # L85: 1153.5238 -> 1188.1052
$ws.Cells.Item(85, 13).Value = 394.6667  # M85: 392 -> 394.6667
$ws.Cells.Item(85, 14).Value = -3684.1052  # N85: -3649.5238 -> -3684.1052

# Row 132 (LTW)
$ws.Cells.Item(132, 8).Value = 4973.5415  # H132: 4392 -> 4973.5415
$ws.Cells.Item(132, 9).Value = 4986.0586  # I132: 4120.2607 -> 4986.0586
$ws.Cells.Item(132, 10).Value = 4943.143  # J132: 5433.6665 -> 4943.143
$ws.Cells.Item(132, 11).Value = 14958.1758  # K132: 12360.7821 -> 14958.1758
$ws.Cells.Item(132, 12).Value = 14829.429  # L132: 16300.9995 -> 14829.429
$ws.Cells.Item(132, 13).Value = -12428.1758  # M132: -9830.7821 -> -12428.1758
$ws.Cells.Item(132, 14).Value = -19889.429  # N132: -21360.9995 -> -19889.429

$ws = $wb.Worksheets.Item("WVR")
# Row 100 (WVR)
$ws.Cells.Item(100, 8).Value = 3209.8462  # H100: 1413.2307 -> 3209.8462
$ws.Cells.Item(100, 9).Value = 3438.8572  # I100: 1252.1111 -> 3438.8572
$ws.Cells.Item(100, 10).Value = 2942.6667  # J100: 1775.75 -> 2942.6667
$ws.Cells.Item(100, 11).Value = 2942.6667  # K100: 2504.2222 -> 2942.6667
$ws.Cells.Item(100, 12).Value = 5885.3334  # L100: 3551.5 -> 5885.3334
$ws.Cells.Item(100, 13).Value = -6336.7144  # M100: -1963.2222 -> -6336.7144
$ws.Cells.Item(100, 14).Value = -6967.3334  # N100: -4633.5 -> -6967.3334
